$d = $word.ActiveDocument

# --- Step 1: split the run "spending unnecessary time and effort in  determining"
# into two pieces at the boundary right before "in  determining", using a
# formatting no-op (Bold on, Bold off) to force Word to materialize a run break
# without altering the visible formatting.
$r = $d.Content
$r.Find.Execute("spending unnecessary time and effort in  determining")
$start = $r.Start
$end = $r.End
$splitLen = "spending unnecessary time and effort ".Length
$splitPos = $start + $splitLen
$rngTail = $d.Range($splitPos, $end)
$rngTail.Font.Bold = 1
$rngTail.Font.Bold = 0

# --- Step 2: normalize the double space ("in  determining" -> "in determining")
$r2 = $d.Content
$r2.Find.Execute("in  determining")
$spacePos = $r2.Start + 2
$rngSpace = $d.Range($spacePos, $spacePos + 1)
$rngSpace.Delete()

# --- Step 3: the text-delete above re-merges various same-formatted runs in
# this paragraph together (standard Word run-consolidation side effect), so
# restore each run break that must survive using the same formatting no-op
# trick.

# 3a: restore the break between "...effort " and "in determining"
$r3 = $d.Content
$r3.Find.Execute("spending unnecessary time and effort in determining")
$start3 = $r3.Start
$end3 = $r3.End
$splitLen3 = "spending unnecessary time and effort ".Length
$splitPos3 = $start3 + $splitLen3
$rngTail3 = $d.Range($splitPos3, $end3)
$rngTail3.Font.Bold = 1
$rngTail3.Font.Bold = 0

# 3b: restore the break between " which " and "movie to watch.  "
$r4 = $d.Content
$r4.Find.Execute("which movie to watch")
$start4 = $r4.Start
$splitLen4 = "which ".Length
$splitPos4 = $start4 + $splitLen4
$rngTail4 = $d.Range($splitPos4, $splitPos4 + "movie to watch.  ".Length)
$rngTail4.Font.Bold = 1
$rngTail4.Font.Bold = 0
